# Config file update: rework the "ScrapyTargetFolder" section on the
# Settings sheet.
#   - B11 (ScrapyTargetFolder's value) now points at a dedicated
#     "DataFromDevices" sub-folder, and gets a real description in C11.
#   - Two new name/value/description rows are appended:
#       row 12: QueueForProcessFolder
#       row 13: QueueForDeviceList

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Row 11 : ScrapyTargetFolder -------------------------------------------------
$ws.Range("B11").Value = 'C:\Temp\Zhengxin\DataFromDevices\'
$ws.Range("C11").Value = '直接从设备下载，未加工的文件存储位置，每次下载都会创建一个新的文件夹'
$ws.Range("C11").Font.Name = "宋体"

# --- Row 12 : QueueForProcessFolder (new) -----------------------------------------
$ws.Range("A12").Value = "QueueForProcessFolder"
$ws.Range("B12").Value = "QueueForProcessFolder"
$ws.Range("C12").Value = "下载完毕所在文件夹（每次下载都会创建子文件夹）"
$ws.Range("C12").Font.Name = "宋体"

# --- Row 13 : QueueForDeviceList (new) --------------------------------------------
$ws.Range("A13").Value = "QueueForDeviceList"
$ws.Range("B13").Value = "QueueForDeviceList"
$ws.Range("C13").Value = "待读取设备队列，再下载的流程Zhengxin-Download里面使用"
$ws.Range("C13").Font.Name = "宋体"

# Match the author's final selection on the Settings sheet.
$ws.Range("C21").Select()
